$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Amount = 1 and Tested? = "Y" for the lumber/clay/ore/stone entity-effect test rows (30-32)
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = "Y"

$ws.Range("C31").Value = 1
$ws.Range("D31").Value = "Y"

$ws.Range("C32").Value = 1
$ws.Range("D32").Value = "Y"

# Set Amount = 1 only (Tested? stays empty) for rows 33-35
$ws.Range("C33").Value = 1
$ws.Range("C34").Value = 1
$ws.Range("C35").Value = 1

# Move the active selection to A36
$ws.Range("A36").Select()
